$wb = $excel.ActiveWorkbook

# Sheet 1: "展览" (index 1) - update column F values for several rows
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 1968
$ws1.Range("F4").Value = 119
$ws1.Range("F6").Value = 15
$ws1.Range("F7").Value = 1644
$ws1.Range("F8").Value = 23
$ws1.Range("F9").Value = 652
$ws1.Range("F12").Value = 22
$ws1.Range("F14").Value = 225
$ws1.Range("F15").Value = 8
$ws1.Range("F18").Value = 131
$ws1.Range("F19").Value = 3818
$ws1.Range("F22").Value = 433
$ws1.Range("F23").Value = 348
$ws1.Range("F24").Value = 706
$ws1.Range("F25").Value = 455
$ws1.Range("F28").Value = 1610
$ws1.Range("F29").Value = 19
$ws1.Range("F30").Value = 155
$ws1.Range("F31").Value = 3

# Sheet 4: "全部类型" (index 4) - update column F values (same rows, plus F27)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 1968
$ws4.Range("F4").Value = 119
$ws4.Range("F6").Value = 15
$ws4.Range("F7").Value = 1644
$ws4.Range("F8").Value = 23
$ws4.Range("F9").Value = 652
$ws4.Range("F12").Value = 22
$ws4.Range("F14").Value = 225
$ws4.Range("F15").Value = 8
$ws4.Range("F18").Value = 131
$ws4.Range("F19").Value = 3818
$ws4.Range("F22").Value = 433
$ws4.Range("F23").Value = 348
$ws4.Range("F24").Value = 706
$ws4.Range("F25").Value = 455
$ws4.Range("F27").Value = 30
$ws4.Range("F28").Value = 1610
$ws4.Range("F29").Value = 19
$ws4.Range("F30").Value = 155
$ws4.Range("F31").Value = 3
